$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 36844.41
$ws.Range("I98").Value = 968.1818
$ws.Range("J98").Value = 72720.63
$ws.Range("K98").Value = 968.1818
$ws.Range("L98").Value = 72720.63
$ws.Range("M98").Value = 529.8182
$ws.Range("N98").Value = -75716.63
# Row 117
$ws.Range("H117").Value = 37709.332
$ws.Range("J117").Value = 37709.332
$ws.Range("L117").Value = 37709.332
$ws.Range("N117").Value = -46887.332
# Row 120
$ws.Range("H120").Value = 49714
$ws.Range("J120").Value = 49714
$ws.Range("L120").Value = 49714
$ws.Range("N120").Value = -59390
# Row 122
$ws.Range("H122").Value = 36844.41
$ws.Range("I122").Value = 968.1818
$ws.Range("J122").Value = 72720.63
$ws.Range("K122").Value = 2904.5454
$ws.Range("L122").Value = 218161.89
$ws.Range("M122").Value = -454.5454
$ws.Range("N122").Value = -223061.89
# Row 130
$ws.Range("H130").Value = 46004.8
$ws.Range("J130").Value = 46004.8
$ws.Range("L130").Value = 46004.8
$ws.Range("N130").Value = -56044.8

$ws = $wb.Worksheets.Item("ARM")
# Row 80
$ws.Range("H80").Value = 50609.43
$ws.Range("J80").Value = 50609.43
$ws.Range("L80").Value = 50609.43
$ws.Range("N80").Value = -52605.43
# Row 83
$ws.Range("H83").Value = 50609.43
$ws.Range("J83").Value = 50609.43
$ws.Range("L83").Value = 151828.29
$ws.Range("N83").Value = -161812.29
# Row 107
$ws.Range("H107").Value = 34327.332
$ws.Range("J107").Value = 34327.332
$ws.Range("L107").Value = 34327.332
$ws.Range("N107").Value = -42007.332
# Row 109
$ws.Range("H109").Value = 40120.332
$ws.Range("J109").Value = 40120.332
$ws.Range("L109").Value = 40120.332
$ws.Range("N109").Value = -42894.332
# Row 110
$ws.Range("H110").Value = 1637.1305
$ws.Range("I110").Value = 1591.2222
$ws.Range("J110").Value = 1802.4
$ws.Range("K110").Value = 1591.2222
$ws.Range("L110").Value = 1802.4
$ws.Range("M110").Value = 453.7778000000001
$ws.Range("N110").Value = -5892.4
# Row 117
$ws.Range("H117").Value = 43816
$ws.Range("J117").Value = 43816
$ws.Range("L117").Value = 43816
$ws.Range("N117").Value = -52994
# Row 118
$ws.Range("H118").Value = 25000
$ws.Range("J118").Value = 25000
$ws.Range("L118").Value = 25000
$ws.Range("N118").Value = -28314
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 125
$ws.Range("H125").Value = 50353.5
$ws.Range("J125").Value = 50353.5
$ws.Range("L125").Value = 50353.5
$ws.Range("N125").Value = -60193.5
# Row 130
$ws.Range("H130").Value = 41248
$ws.Range("J130").Value = 41248
$ws.Range("L130").Value = 41248
$ws.Range("N130").Value = -51288
# Row 131
$ws.Range("H131").Value = 51617
$ws.Range("J131").Value = 51617
$ws.Range("L131").Value = 51617
$ws.Range("N131").Value = -61697
# Row 132
$ws.Range("H132").Value = 21740716
$ws.Range("I132").Value = 27778838
$ws.Range("J132").Value = 3479.2
$ws.Range("K132").Value = 83336514
$ws.Range("L132").Value = 10437.6
$ws.Range("M132").Value = -83333984
$ws.Range("N132").Value = -15497.6

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1400.7273
$ws.Range("I94").Value = 1167.5555
$ws.Range("K94").Value = 1167.5555
$ws.Range("M94").Value = -716.5554999999999
# Row 117
$ws.Range("H117").Value = 44999
$ws.Range("J117").Value = 44999
$ws.Range("L117").Value = 44999
$ws.Range("N117").Value = -54177
# Row 125
$ws.Range("H125").Value = 50772
$ws.Range("J125").Value = 50772
$ws.Range("L125").Value = 50772
$ws.Range("N125").Value = -60612
# Row 126
$ws.Range("H126").Value = 44002.668
$ws.Range("J126").Value = 44002.668
$ws.Range("L126").Value = 44002.668
$ws.Range("N126").Value = -53882.668
# Row 130
$ws.Range("H130").Value = 50779.5
$ws.Range("J130").Value = 50779.5
$ws.Range("L130").Value = 50779.5
$ws.Range("N130").Value = -60819.5

$ws = $wb.Worksheets.Item("CRP")
# Row 116
$ws.Range("H116").Value = 42362.25
$ws.Range("J116").Value = 42362.25
$ws.Range("L116").Value = 42362.25
$ws.Range("N116").Value = -51540.25
# Row 132
$ws.Range("H132").Value = 57201.152
$ws.Range("I132").Value = 2296.2
$ws.Range("J132").Value = 132071.55
$ws.Range("K132").Value = 6888.599999999999
$ws.Range("L132").Value = 396214.65
$ws.Range("M132").Value = -4358.599999999999
$ws.Range("N132").Value = -401274.65

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 4128.4414
$ws.Range("I5").Value = 20421.4
$ws.Range("K5").Value = 61264.2
$ws.Range("M5").Value = -61152.2
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 133
$ws.Range("H133").Value = 7469.5454
$ws.Range("J133").Value = 6923.5713
$ws.Range("L133").Value = 20770.7139
$ws.Range("N133").Value = -30890.7139
# Row 135
$ws.Range("H135").Value = 4128.4414
$ws.Range("I135").Value = 20421.4
$ws.Range("K135").Value = 183792.6
$ws.Range("M135").Value = -181257.6

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 7755
$ws.Range("I126").Value = 13281.667
$ws.Range("J126").Value = 1537.5
$ws.Range("K126").Value = 39845.001
$ws.Range("L126").Value = 4612.5
$ws.Range("M126").Value = -37375.001
$ws.Range("N126").Value = -9552.5
# Row 130
$ws.Range("H130").Value = 45985.125
$ws.Range("J130").Value = 45985.125
$ws.Range("L130").Value = 45985.125
$ws.Range("N130").Value = -56025.125
# Row 138
$ws.Range("H138").Value = 51545.453
$ws.Range("J138").Value = 51545.453
$ws.Range("L138").Value = 51545.453
$ws.Range("N138").Value = -61825.453

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3573.2666
$ws.Range("I132").Value = 1380.4
$ws.Range("J132").Value = 4669.7
$ws.Range("K132").Value = 4141.200000000001
$ws.Range("L132").Value = 14009.1
$ws.Range("M132").Value = -1611.200000000001
$ws.Range("N132").Value = -19069.1
# Row 136
$ws.Range("H136").Value = 2153.2
$ws.Range("I136").Value = 1467
$ws.Range("J136").Value = 4898
$ws.Range("K136").Value = 4401
$ws.Range("L136").Value = 14694
$ws.Range("M136").Value = -1851
$ws.Range("N136").Value = -19794

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1950.2106
$ws.Range("I132").Value = 1262.5834
$ws.Range("J132").Value = 3129
$ws.Range("K132").Value = 3787.7502
$ws.Range("L132").Value = 9387
$ws.Range("M132").Value = -1257.7502
$ws.Range("N132").Value = -14447
# Row 133
$ws.Range("H133").Value = 91354.664
$ws.Range("J133").Value = 91354.664
$ws.Range("L133").Value = 91354.664
$ws.Range("N133").Value = -101474.664
# Row 136
$ws.Range("H136").Value = 19318.91
$ws.Range("I136").Value = 48412.57
$ws.Range("K136").Value = 145237.71
$ws.Range("M136").Value = -142687.71
